# Update status of saved-query tracking item #27 (row 28) to "Complete",
# record its completion date, and clear the AutoFilter so all previously
# filtered-out rows (Status = In Progress / Not Started) become visible
# again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Mark item #27 complete and stamp the completion date ------------------
$ws.Range("F28").Value = "Complete"

$ws.Range("H28").NumberFormat = "mm-dd-yy"
$ws.Range("H28").Value = $ws.Range("G28").Value2

# --- Clear the AutoFilter criteria (Status = In Progress / Not Started) ----
# so every row is shown again (this also un-hides the rows Excel hid to
# implement the filter).
$ws.AutoFilterMode = $false
$ws.Range("A1:H28").AutoFilter()

# --- Update the saved cursor/selection position -----------------------------
$ws.Range("F24").Select()
